$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-08-06 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-08-07 Thursday", 2)
$d.Content.Find.Execute("43÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "11÷9=", 2)
$d.Content.Find.Execute("66÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "87÷3=", 2)
$d.Content.Find.Execute("96÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "14÷3=", 2)
$d.Content.Find.Execute("69÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "12÷9=", 2)
$d.Content.Find.Execute("23÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "50÷5=", 2)
$d.Content.Find.Execute("13÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "68÷6=", 2)
$d.Content.Find.Execute("10÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "18÷8=", 2)
$d.Content.Find.Execute("26÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "74÷7=", 2)
$d.Content.Find.Execute("32÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "56÷4=", 2)
$d.Content.Find.Execute("84÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "27÷8=", 2)
$d.Content.Find.Execute("56÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "92÷5=", 2)
$d.Content.Find.Execute("53÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "76÷9=", 2)
$d.Content.Find.Execute("40÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "67÷9=", 2)
$d.Content.Find.Execute("63÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "58÷9=", 2)
$d.Content.Find.Execute("46÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "45÷4=", 2)
$d.Content.Find.Execute("66÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "11÷7=", 2)
$d.Content.Find.Execute("27÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "32÷3=", 2)
$d.Content.Find.Execute("14÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "51÷4=", 2)
$d.Content.Find.Execute("25÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "56÷2=", 2)
$d.Content.Find.Execute("64÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "90÷8=", 2)
$d.Content.Find.Execute("58÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "42÷2=", 2)
$d.Content.Find.Execute("66÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "51÷6=", 2)
$d.Content.Find.Execute("28÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "79÷6=", 2)
$d.Content.Find.Execute("53÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "19÷6=", 2)
$d.Content.Find.Execute("84÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "73÷8=", 2)
